# Review_446.docx edit: update paper review content to the new paper
# "Multimodal Latent Language Modeling with Next-Token Diffusion" (26.04.25).
$d = $word.ActiveDocument

# ---- New text content (verbatim from the target revision) ----
$DATE_TEXT = ' המאמר היומי של מייק: 26.04.25'
$TITLE_TEXT = 'Multimodal Latent Language Modeling with Next-Token Diffusion'
$P2_TEXT = 'היום שבת והסקירה של היום תהיה קלילה ודי קצרה. הסקירה תתמקד במודלים מולטי-מודליים גנרטיביים המסוגלים "להבין״ וליצור דאטה מכמה מודליות כלומר טקסט, תמונות, אודיו וכדומה. המאמר למעשה בעצם משדך מודלים לטנטיים גנרטיבים עבור דאטה טקסטואלי ועבור דאטה רציף יותר (למרות שגם הוא discretized). המחברים עושים זאת באמצעות אימון של מודלי דיפוזיה גנרטיבים עבור סוגי דאטה שונים במרחב הלטנטי. כלומר המודל מאומן לגנרט ייצוגים לטנטיים עבור דאטה טקסטואלי ועבור דאטה כמו אודיו ותמונות.'
$P3_TEXT = 'להבדיל ממאמרים רבים אחרים המחברים מאמנים לא רק את המודל הגנרטיבי המולטימודלי אלא מאמנים גם מודל אמבדינג להפקה של ייצוגים לטנטיים של דאטה ממודליות שונות. בדרך כלל מודל האמבדינג במודלי דיפוזיה הוא מסוג VAE (שזה Variational Autoencoder) והמחברים מציעים מודיפיקציה קלה ל-VAE. במקום שהאנקודר (הקלט אליו הוא דאטה) של VAE יגנרט את וקטורי התוחלות השונויות של הווקטור הלטנטי הוא מגנרט רק וקטור התוחלות כאשר השוניות מוגרלות ההתפלגות גאוסית עם שונות נתונה (הייפרפרמטר). לדעת המחברים זה מונע קריסה(איפוס) של וקטור השונויות הנוצר על ידי האנקודר שפוגע בגיוון התמונות שהמודל מגנרט. '
$P4_TEXT = 'המחברים מאמנים VAE עבור דאטה לא טקסטואלי בלבד. תמונה או אודיו מחולקת לטוקנים (פאצ''ים לתמונות ולמקטעים בזמן לאודיו) ומוזנים למודל כדאטה סדרתית. שימו לב המודל מסתכל על דאטה בכל מודליות כמו דאטה סדרתי. זה מאוד טריוויאלי לדאטה טקסטואלי ולאודיו כי יש שם סדר אינהרנטי ברור. בתמונות גם יש סדר אבל הוא יכול לבוא בכמה צורות: כלומר ניתן לתאר תמונה כסדרה של פאצ''ים בכמה צורות (למשל משמאל לימין ולמעלה למטה וגם מימין לשמאל ומלטה למעלה). '
$P5_TEXT = 'מודל דיפוזיה לדאטה לא טקסטואלי מאומן לנקות את הרעש מהדאטה (denoising) בהינתן הייצוג הלטנטי שלו(המורעש) ושל ההקשר (כל המודלים במאמר כמובן אוטורגרסיביים). לאחר מכן הוקטור הלטנטי הנקי מוזן לדקודר של VAE לשחזור הדאטה כאשר המטרה של המודל המאומן היא לשחזר את הדאטה כמה שיותר טוב. עבור דאטה טקסטואלי ההרעשה מופעלת על האמבדינג של הטוקנים הטקסטואליים ומודל דיפוזיה מאומן לשחזר אותם. עבור דאטה טקסטואלי מאמנים עוד שכבה לינארית שמטרתה למפות את הוקטור הלטנטי למרחב הטוקנים הטקסטואליים (סופטמקס בגודל של מילון).  דרך אגב מודלי דיפוזיה מאומנים יחד עם ה-VAE (אנקודר ודקודר). '
$P6_TEXT = 'כדי להפריד בין דאטה טקסטואלי ולא טקסטואלי המחברים מאמנים טוקנים המפרידים בין דאטה השייך למודליות שונות.'
$URL_TEXT = 'https://arxiv.org/abs/2412.08635'

# ---- Paragraph 1: date line + title line (two <w:t> runs separated by <w:br/>) ----
# Only the day part of the date changes (28 -> 26) and the paper title changes entirely.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("28.04.25", $false, $false, $false, $false, $false, $true, 1, $false, "26.04.25", 2) | Out-Null
$find.Execute("Around the World in 80 Timesteps: A Generative Approach to Global Visual Geolocation", $false, $false, $false, $false, $false, $true, 1, $false, $TITLE_TEXT, 2) | Out-Null
# Sanity-check the full first line now reads as expected (date run + title run).
if ($d.Paragraphs.Item(1).Range.Text -notmatch [regex]::Escape($DATE_TEXT.Trim())) {
    throw "Paragraph 1 date text not updated as expected"
}

# ---- Paragraphs 2-6: replace each paragraph's body text in place ----
function Set-ParagraphText($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

Set-ParagraphText 2 $P2_TEXT
Set-ParagraphText 3 $P3_TEXT
Set-ParagraphText 4 $P4_TEXT
Set-ParagraphText 5 $P5_TEXT
Set-ParagraphText 6 $P6_TEXT

# ---- Remove the three paragraphs that followed (old paragraphs 7-9) ----
$delStart = $d.Paragraphs.Item(7).Range.Start
$delEnd = $d.Paragraphs.Item(9).Range.End
$d.Range($delStart, $delEnd).Delete() | Out-Null

# ---- Final paragraph: replace the arxiv link text (now paragraph 7) ----
Set-ParagraphText 7 $URL_TEXT

# ---- Final sanity checks ----
if ($d.Paragraphs.Count -ne 7) {
    throw "Expected 7 paragraphs after edit, found $($d.Paragraphs.Count)"
}
if ($d.Paragraphs.Item(7).Range.Text.Trim() -ne $URL_TEXT) {
    throw "Final paragraph does not contain the expected link text"
}

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
